$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5077
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 9999
$ws.Range("K18").Value = 155
$ws.Range("L18").Value = 9999
$ws.Range("M18").Value = 129
$ws.Range("N18").Value = -10567
$ws.Range("H19").Value = 1074.6666
$ws.Range("I19").Value = 1033.4
$ws.Range("J19").Value = 1104.1428
$ws.Range("K19").Value = 1033.4
$ws.Range("L19").Value = 1104.1428
$ws.Range("M19").Value = -858.4000000000001
$ws.Range("N19").Value = -1454.1428
$ws.Range("H38").Value = 59
$ws.Range("I38").Value = 59
$ws.Range("K38").Value = 177
$ws.Range("M38").Value = 195
$ws.Range("H41").Value = 380.84616
$ws.Range("J41").Value = 2999
$ws.Range("L41").Value = 2999
$ws.Range("N41").Value = -3879
$ws.Range("H42").Value = 1305.8889
$ws.Range("I42").Value = 36.285713
$ws.Range("J42").Value = 5749.5
$ws.Range("K42").Value = 108.857139
$ws.Range("L42").Value = 17248.5
$ws.Range("M42").Value = 121.142861
$ws.Range("N42").Value = -17708.5
$ws.Range("H43").Value = 2679.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2679.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2679.8
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2817.8
$ws.Range("H62").Value = 2524
$ws.Range("I62").Value = 2533
$ws.Range("K62").Value = 2533
$ws.Range("M62").Value = -1909
$ws.Range("H65").Value = 2524
$ws.Range("I65").Value = 2533
$ws.Range("K65").Value = 12665
$ws.Range("M65").Value = -9545
$ws.Range("H93").Value = 40475
$ws.Range("J93").Value = 40475
$ws.Range("L93").Value = 40475
$ws.Range("N93").Value = -45467
$ws.Range("H94").Value = 288
$ws.Range("I94").Value = 288
$ws.Range("K94").Value = 288
$ws.Range("M94").Value = 163
$ws.Range("H106").Value = 13818.895
$ws.Range("I106").Value = 15221.117
$ws.Range("J106").Value = 1900
$ws.Range("K106").Value = 15221.117
$ws.Range("L106").Value = 1900
$ws.Range("M106").Value = -14590.117
$ws.Range("N106").Value = -3162
$ws.Range("H107").Value = 3605.36
$ws.Range("I107").Value = 3717.85
$ws.Range("J107").Value = 3155.4
$ws.Range("K107").Value = 3717.85
$ws.Range("L107").Value = 3155.4
$ws.Range("M107").Value = -1797.85
$ws.Range("N107").Value = -6995.4
$ws.Range("H116").Value = 19853.262
$ws.Range("I116").Value = 29101
$ws.Range("K116").Value = 29101
$ws.Range("M116").Value = -25659
$ws.Range("H138").Value = 3016
$ws.Range("J138").Value = 4133.3687
$ws.Range("L138").Value = 12400.1061
$ws.Range("N138").Value = -22680.1061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3254
$ws.Range("I21").Value = 3254
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3254
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -2880
$ws.Range("N21").ClearContents()
$ws.Range("H102").Value = 2797.6667
$ws.Range("I102").Value = 2797.6667
$ws.Range("K102").Value = 2797.6667
$ws.Range("M102").Value = -1175.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1177290
$ws.Range("I22").Value = 602.1429000000001
$ws.Range("K22").Value = 602.1429000000001
$ws.Range("M22").Value = -429.1429000000001
$ws.Range("H105").Value = 860.6667
$ws.Range("I105").Value = 860.6667
$ws.Range("K105").Value = 860.6667
$ws.Range("M105").Value = 886.3333
$ws.Range("H134").Value = 6350.6343
$ws.Range("I134").Value = 6129.2583
$ws.Range("K134").Value = 18387.7749
$ws.Range("M134").Value = -15852.7749

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 534.73914
$ws.Range("I22").Value = 340.125
$ws.Range("J22").Value = 638.5333000000001
$ws.Range("K22").Value = 340.125
$ws.Range("L22").Value = 638.5333000000001
$ws.Range("M22").Value = 9.875
$ws.Range("N22").Value = -1338.5333
$ws.Range("H58").Value = 3458.2068
$ws.Range("I58").Value = 2616.4
$ws.Range("J58").Value = 5328.8887
$ws.Range("K58").Value = 2616.4
$ws.Range("L58").Value = 5328.8887
$ws.Range("M58").Value = -2413.4
$ws.Range("N58").Value = -5734.8887
$ws.Range("H103").Value = 42499.75
$ws.Range("I103").Value = 17500
$ws.Range("K103").Value = 17500
$ws.Range("M103").Value = -16328
$ws.Range("H136").Value = 3458.2068
$ws.Range("I136").Value = 2616.4
$ws.Range("J136").Value = 5328.8887
$ws.Range("K136").Value = 7849.200000000001
$ws.Range("L136").Value = 15986.6661
$ws.Range("M136").Value = -5299.200000000001
$ws.Range("N136").Value = -21086.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 451.1111
$ws.Range("I14").Value = 451.1111
$ws.Range("K14").Value = 1353.3333
$ws.Range("M14").Value = -1180.3333
$ws.Range("H45").Value = 315
$ws.Range("J45").Value = 315
$ws.Range("L45").Value = 945
$ws.Range("N45").Value = -2009
$ws.Range("H107").Value = 1500.5714
$ws.Range("I107").Value = 756
$ws.Range("K107").Value = 2268
$ws.Range("M107").Value = -348
$ws.Range("H121").Value = 3394.2727
$ws.Range("J121").Value = 3413.8823
$ws.Range("L121").Value = 10241.6469
$ws.Range("N121").Value = -12861.6469
$ws.Range("H133").Value = 22559.273
$ws.Range("I133").Value = 6615.2856
$ws.Range("K133").Value = 19845.8568
$ws.Range("M133").Value = -14785.8568
$ws.Range("H140").Value = 10421867
$ws.Range("I140").Value = 14708018
$ws.Range("K140").Value = 44124054
$ws.Range("M140").Value = -44118874

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2493.7368
$ws.Range("J122").Value = 3492.8
$ws.Range("L122").Value = 10478.4
$ws.Range("N122").Value = -15378.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 251438.25
$ws.Range("I7").Value = 334334.34
$ws.Range("K7").Value = 334334.34
$ws.Range("M7").Value = -334222.34
$ws.Range("H22").Value = 1258
$ws.Range("J22").Value = 1377.3334
$ws.Range("L22").Value = 1377.3334
$ws.Range("N22").Value = -1967.3334
$ws.Range("H27").Value = 1258
$ws.Range("J27").Value = 1377.3334
$ws.Range("L27").Value = 1377.3334
$ws.Range("N27").Value = -1591.3334
$ws.Range("H40").Value = 2567.6667
$ws.Range("I40").Value = 2202
$ws.Range("K40").Value = 2202
$ws.Range("M40").Value = -2066
$ws.Range("H93").Value = 1233.8572
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752
$ws.Range("H100").Value = 2249.75
$ws.Range("I100").Value = 1999.5
$ws.Range("K100").Value = 1999.5
$ws.Range("M100").Value = -1458.5
$ws.Range("H126").Value = 251438.25
$ws.Range("I126").Value = 334334.34
$ws.Range("K126").Value = 1003003.02
$ws.Range("M126").Value = -1000533.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 26199.75
$ws.Range("J95").Value = 26199.75
$ws.Range("L95").Value = 26199.75
$ws.Range("N95").Value = -31691.75
$ws.Range("H100").Value = 2818.4
$ws.Range("I100").Value = 1366.3334
$ws.Range("K100").Value = 2732.6668
$ws.Range("M100").Value = -2191.6668
$ws.Range("H107").Value = 621.7917
$ws.Range("J107").Value = 879.6667
$ws.Range("L107").Value = 2639.0001
$ws.Range("N107").Value = -6479.0001
$ws.Range("H113").Value = 422
$ws.Range("I113").Value = 414.14285
$ws.Range("K113").Value = 1242.42855
$ws.Range("M113").Value = 927.5714499999999
$ws.Range("H126").Value = 1782.25
$ws.Range("I126").Value = 1498.35
$ws.Range("K126").Value = 4495.049999999999
$ws.Range("M126").Value = -2025.049999999999
$ws.Range("H132").Value = 2859.195
$ws.Range("I132").Value = 2626.697
$ws.Range("J132").Value = 3818.25
$ws.Range("K132").Value = 7880.091
$ws.Range("L132").Value = 11454.75
$ws.Range("M132").Value = -5350.091
$ws.Range("N132").Value = -16514.75
